# Applies the hashcode metadata update to the "hashcode.csv" sheet.
# Each cell in column B holds an MD5-like hash string; this commit
# regenerates a number of those hash values (the cell in column A,
# holding the associated code, stays unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$updates = @{
    "B154" = "6b15316edc1cc092b4abac42be90bd28"
    "B160" = "a971ea9eb8c3823f3586968e3793190b"
    "B169" = "4da83de0fa8baa0c3e34ef948fa497bf"
    "B222" = "b2c2d7b0c6e1e482e2baebfaa3e80238"
    "B227" = "811e4b110a2cffba77fce045c7017d73"
    "B229" = "67e8de9238b1d980854c534789e8446c"
    "B232" = "869c621bbced2dd1e9009bcaac137d49"
    "B284" = "afc91a4d0896544a39504d970bebe301"
    "B423" = "0841f66eec1f7caf51680bed6f5054c6"
    "B486" = "7c7e26fef28b133513b0e1d817db11ed"
    "B526" = "46abcc7d85f2732d753478da077c6dad"
    "B545" = "caed40e30b8d326c9ee29159f49801d9"
    "B565" = "6dae6fa19d878e3e786208dc34f13627"
    "B578" = "c2773ef09b571a4d55e3f514b1138e7d"
    "B584" = "90e9978e5fac4cdc1c413f6cc4049a3c"
    "B692" = "87f7d8c8d5f14748512c9245c79f6ea6"
    "B697" = "e992428de39ad6cc52cb72f089587295"
    "B712" = "c73244e4d02da93b2f5418460dd36c9d"
    "B715" = "d174fa8fbca0c777f41402c2571309ad"
    "B823" = "d05f60cb7fe7ed68b218c83ac767a514"
    "B827" = "828dfcdbe017b46b27ba6a91372baea2"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
